$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 106 ("Plomo") entirely; this shifts rows 107:120 up by one,
# matching the diff which removes the "Plomo" entry and shifts every
# subsequent song's Title/Lyrics up by one row (row 120 disappears).
$ws.Rows.Item(106).Delete()
